$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "63.428.72"
$ws.Cells.Item(2,5).Value = "  +1.36%  "
$ws.Cells.Item(3,4).Value = "2.542.50"
$ws.Cells.Item(3,5).Value = "  +4.55%  "
$ws.Cells.Item(4,5).Value = "  -0.03%  "
$ws.Cells.Item(5,4).Value = "'568.32"
$ws.Cells.Item(5,5).Value = "  +1.80%  "
$ws.Cells.Item(6,4).Value = "'150.53"
$ws.Cells.Item(6,5).Value = "  +7.98%  "
$ws.Cells.Item(7,5).Value = "  -0.01%  "
$ws.Cells.Item(8,5).Value = "  -0.01%  "
$ws.Cells.Item(9,4).Value = "2.541.35"
$ws.Cells.Item(9,5).Value = "  +4.54%  "
$ws.Cells.Item(10,5).Value = "  +1.39%  "
$ws.Cells.Item(11,5).Value = "  -1.79%  "
$ws.Cells.Item(12,5).Value = "  +1.10%  "
$ws.Cells.Item(13,4).Value = "'0.357"
$ws.Cells.Item(14,5).Value = "  +7.41%  "
$ws.Cells.Item(15,4).Value = "2.995.60"
$ws.Cells.Item(15,5).Value = "  +4.56%  "
$ws.Cells.Item(16,4).Value = "63.323.32"
$ws.Cells.Item(16,5).Value = "  +1.37%  "
$ws.Cells.Item(17,5).Value = "  +1.34%  "
$ws.Cells.Item(18,4).Value = "2.521.19"
$ws.Cells.Item(18,5).Value = "  +3.68%  "
$ws.Cells.Item(19,4).Value = "'11.59"
$ws.Cells.Item(19,5).Value = "  +2.92%  "
$ws.Cells.Item(20,4).Value = "'339.45"
$ws.Cells.Item(20,5).Value = "  -2.21%  "
$ws.Cells.Item(21,5).Value = "  +2.72%  "
$ws.Cells.Item(22,5).Value = "  -0.85%  "
$ws.Cells.Item(23,5).Value = "  +0.08%  "
$ws.Cells.Item(24,4).Value = "'65.90"
$ws.Cells.Item(24,5).Value = "  +0.50%  "
$ws.Cells.Item(25,5).Value = "  -2.03%  "
$ws.Cells.Item(26,4).Value = "'1.56"
$ws.Cells.Item(26,5).Value = "  +14.62%  "
$ws.Cells.Item(27,4).Value = "'1.59"
$ws.Cells.Item(27,5).Value = "  +1.63%  "
$ws.Cells.Item(28,5).Value = "  +0.16%  "
$ws.Cells.Item(29,4).Value = "'8.44"
$ws.Cells.Item(29,5).Value = "  +3.52%  "
$ws.Cells.Item(30,5).Value = "  +10.00%  "
$ws.Cells.Item(31,4).Value = "0.0₃0817"
$ws.Cells.Item(31,5).Value = "  +4.07%  "
$ws.Cells.Item(32,5).Value = "  +1.87%  "
$ws.Cells.Item(34,4).Value = "'1.57"
$ws.Cells.Item(34,5).Value = "  +8.39%  "
$ws.Cells.Item(35,4).Value = "'419.82"
$ws.Cells.Item(35,5).Value = "  +10.68%  "
$ws.Cells.Item(36,4).Value = "'0.404"
$ws.Cells.Item(36,5).Value = "  +1.57%  "
$ws.Cells.Item(37,5).Value = "  +2.16%  "
$ws.Cells.Item(38,4).Value = "'4.41"
$ws.Cells.Item(38,5).Value = "  -1.12%  "
$ws.Cells.Item(39,5).Value = "  -0.02%  "
$ws.Cells.Item(40,5).Value = "  +3.63%  "
$ws.Cells.Item(41,4).Value = "'1.00"
$ws.Cells.Item(41,5).Value = "  -0.01%  "
$ws.Cells.Item(42,5).Value = "  -0.33%  "
$ws.Cells.Item(43,4).Value = "'153.61"
$ws.Cells.Item(43,5).Value = "  +6.08%  "
$ws.Cells.Item(44,5).Value = "  +2.47%  "
$ws.Cells.Item(45,4).Value = "'20.74"
$ws.Cells.Item(45,5).Value = "  -0.41%  "
$ws.Cells.Item(46,5).Value = "  +2.09%  "
$ws.Cells.Item(47,5).Value = "  +0.87%  "
$ws.Cells.Item(48,4).Value = "'0.0522"
$ws.Cells.Item(48,5).Value = "  +0.13%  "
$ws.Cells.Item(50,4).Value = "'18.43"
$ws.Cells.Item(50,5).Value = "  +3.30%  "
$ws.Cells.Item(51,4).Value = "'1.79"
$ws.Cells.Item(51,5).Value = "  +4.03%  "
